$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (shared string reused by A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 20:40"

# --- Estados Unidos (row 4): refresh counters ---
$ws.Range("B4").Value = 1806255
$ws.Range("C4").Value = 12725
$ws.Range("E4").Value = 1178134
$ws.Range("G4").Value = 517
$ws.Range("H4").Value = 105059

# --- Francia (row 10): refresh counters ---
$ws.Range("B10").Value = 186797
$ws.Range("E10").Value = 90280

# --- Suiza / Sudafrica swap ranking (rows 32-33) ---
# Sudafrica moves up to row 32 with fresh counters; Suiza drops to row 33
# keeping its previous counters unchanged.
$ws.Range("A32").Value = "Sudafrica"
$ws.Range("B32").Value = 30967
$ws.Range("C32").Value = 1727
$ws.Range("D32").Value = 16116
$ws.Range("E32").Value = 14208
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 32
$ws.Range("H32").Value = 643

$ws.Range("A33").Value = "Suiza"
$ws.Range("B33").Value = 30845
$ws.Range("C33").Value = 17
$ws.Range("D33").Value = 28400
$ws.Range("E33").Value = 526
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 1919

# --- Marruecos / Ghana / Malasia swap ranking (rows 63-65) ---
# Marruecos moves up to row 63 with fresh counters; Ghana and Malasia
# each drop one row, keeping their previous counters unchanged.
$ws.Range("A63").Value = "Marruecos"
$ws.Range("B63").Value = 7780
$ws.Range("C63").Value = 66
$ws.Range("D63").Value = 5401
$ws.Range("E63").Value = 2175
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 204

$ws.Range("A64").Value = "Ghana"
$ws.Range("B64").Value = 7768
$ws.Range("C64").Value = 152
$ws.Range("D64").Value = 2540
$ws.Range("E64").Value = 5193
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 35

$ws.Range("A65").Value = "Malasia"
$ws.Range("B65").Value = 7762
$ws.Range("C65").Value = 30
$ws.Range("D65").Value = 6330
$ws.Range("E65").Value = 1317
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 115

# --- Republica de Chipre (row 121): refresh counters ---
$ws.Range("D121").Value = 790
$ws.Range("E121").Value = 137
